# "add drop item list record"
#
# Record_PosList row 1 held the field-name headers and row 2 held the
# X / Y / Z / StayTime pos-list fields; this edit repurposes row 1 as the
# X/Y/Z/StayTime header (so the struct describes a position-list record)
# and pushes the generic "float" placeholders down into row 2 as the
# first data row. It also re-targets the 4 column comments to the new
# header row and switches the active sheet/selection to Record_PosList.

$wb = $excel.ActiveWorkbook
$wsProperty = $wb.Worksheets.Item("Property")
$wsPosList  = $wb.Worksheets.Item("Record_PosList")

# --- swap the L:O header/value rows on Record_PosList -----------------
# Row 1 currently holds the generic "float" placeholders (L1:O1) and
# row 2 holds the real X/Y/Z/StayTime field names (L2:O2); the edit
# swaps them so row 1 becomes the header and row 2 becomes sample data.
$row1Vals = @(
    $wsPosList.Range("L1").Value2,
    $wsPosList.Range("M1").Value2,
    $wsPosList.Range("N1").Value2,
    $wsPosList.Range("O1").Value2
)
$row2Vals = @(
    $wsPosList.Range("L2").Value2,
    $wsPosList.Range("M2").Value2,
    $wsPosList.Range("N2").Value2,
    $wsPosList.Range("O2").Value2
)

$wsPosList.Range("L1").Value = $row2Vals[0]
$wsPosList.Range("M1").Value = $row2Vals[1]
$wsPosList.Range("N1").Value = $row2Vals[2]
$wsPosList.Range("O1").Value = $row2Vals[3]

$wsPosList.Range("L2").Value = $row1Vals[0]
$wsPosList.Range("M2").Value = $row1Vals[1]
$wsPosList.Range("N2").Value = $row1Vals[2]
$wsPosList.Range("O2").Value = $row1Vals[3]

# --- re-target the column comments from row 1 to row 2 ----------------
$wsPosList.Range("L1").Comment.Delete()
$wsPosList.Range("M1").Comment.Delete()
$wsPosList.Range("N1").Comment.Delete()
$wsPosList.Range("O1").Comment.Delete()

$wsPosList.Range("L2").AddComment("强化等级")
$wsPosList.Range("M2").AddComment("强化等级")
$wsPosList.Range("N2").AddComment("强化等级")
$wsPosList.Range("O2").AddComment("镶嵌宝石，逗号分隔")

# --- switch the active sheet / selection -------------------------------
$wsProperty.Range("J33").Select()
$wsPosList.Activate()
$wsPosList.Range("O8").Select()
